$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.443.68'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '2.592.17'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '571.12'
$ws.Range('E5').Value = '  +3.27%  '
$ws.Range('D6').Value = '144.13'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '2.604.40'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('D11').Value = '0.105'
$ws.Range('E11').Value = '  +3.92%  '
$ws.Range('D12').Value = '0.158'
$ws.Range('E12').Value = '  +11.36%  '
$ws.Range('D13').Value = '0.346'
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('D14').Value = '3.052.17'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').Value = '59.401.69'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = '22.73'
$ws.Range('E16').Value = '  +8.78%  '
$ws.Range('E17').Value = '  +4.25%  '
$ws.Range('D18').Value = '2.600.30'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D20').Value = '335.89'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').Value = '10.30'
$ws.Range('E21').Value = '  +1.89%  '
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '64.55'
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('D25').Value = '0.459'
$ws.Range('E25').Value = '  +6.95%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.161'
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.982'
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').Value = '7.32'
$ws.Range('E28').Value = '  +1.97%  '
$ws.Range('D29').Value = '0.0₃0786'
$ws.Range('E29').Value = '  +4.05%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').Value = '6.12'
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('D33').Value = '158.47'
$ws.Range('E33').Value = '  +2.84%  '
$ws.Range('D34').Value = '19.10'
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').Value = '4.08'
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('D36').Value = '1.16'
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('D37').Value = '0.885'
$ws.Range('E37').Value = '  +2.06%  '
$ws.Range('D38').Value = '0.879'
$ws.Range('E38').Value = '  -2.96%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '37.17'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.50'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('D41').Value = '296.66'
$ws.Range('E41').Value = '  +4.56%  '
$ws.Range('D42').Value = '3.70'
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '0.0980'
$ws.Range('E44').Value = '  +2.68%  '
$ws.Range('D45').Value = '0.597'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '19.39'
$ws.Range('E46').Value = '  +3.18%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '0.0540'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').Value = '10.63'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '125.65'
$ws.Range('E49').Value = '  +7.06%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0234'
$ws.Range('E50').Value = '  +2.38%  '
$ws.Range('D51').Value = '18.69'
$ws.Range('E51').Value = '  +3.46%  '
